$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "57.803.58"
Set-TextValue $ws.Range("E2") "  +2.50%  "

Set-TextValue $ws.Range("D3") "3.062.58"

Set-TextValue $ws.Range("E4") "  +0.03%  "

Set-TextValue $ws.Range("D5") "516.39"
Set-TextValue $ws.Range("E5") "  +1.48%  "

Set-TextValue $ws.Range("D6") "142.75"
Set-TextValue $ws.Range("E6") "  +1.93%  "

Set-TextValue $ws.Range("E7") "  +0.05%  "

Set-TextValue $ws.Range("D8") "0.436"
Set-TextValue $ws.Range("E8") "  +1.85%  "

Set-TextValue $ws.Range("D9") "7.27"
Set-TextValue $ws.Range("E9") "  +2.37%  "

Set-TextValue $ws.Range("E10") "  +0.17%  "

Set-TextValue $ws.Range("E11") "  +2.52%  "

Set-TextValue $ws.Range("D12") "3.594.19"
Set-TextValue $ws.Range("E12") "  +2.28%  "

Set-TextValue $ws.Range("E13") "  +3.04%  "

Set-TextValue $ws.Range("E14") "  +3.85%  "

Set-TextValue $ws.Range("E15") "  +1.33%  "

Set-TextValue $ws.Range("D16") "57.846.94"
Set-TextValue $ws.Range("E16") "  +2.65%  "

Set-TextValue $ws.Range("D17") "3.066.35"
Set-TextValue $ws.Range("E17") "  +2.03%  "

Set-TextValue $ws.Range("D18") "6.09"
Set-TextValue $ws.Range("E18") "  +3.20%  "

Set-TextValue $ws.Range("D19") "12.82"
Set-TextValue $ws.Range("E19") "  -0.38%  "

Set-TextValue $ws.Range("D20") "8.09"
Set-TextValue $ws.Range("E20") "  +1.28%  "

Set-TextValue $ws.Range("D21") "332.17"
Set-TextValue $ws.Range("E21") "  +0.43%  "

Set-TextValue $ws.Range("D22") "0.998"
Set-TextValue $ws.Range("E22") "  -0.01%  "

Set-TextValue $ws.Range("E23") "  +0.35%  "

Set-TextValue $ws.Range("D24") "65.40"
Set-TextValue $ws.Range("E24") "  +1.47%  "

Set-TextValue $ws.Range("E25") "  +3.35%  "

Set-TextValue $ws.Range("E26") "  +0.39%  "

Set-TextValue $ws.Range("D27") "0.0₃0902"
Set-TextValue $ws.Range("E27") "  -3.52%  "

Set-TextValue $ws.Range("D28") "6.47"
Set-TextValue $ws.Range("E28") "  +2.20%  "

Set-TextValue $ws.Range("D29") "7.22"
Set-TextValue $ws.Range("E29") "  +5.72%  "

Set-TextValue $ws.Range("E30") "  +1.75%  "

Set-TextValue $ws.Range("D32") "20.67"
Set-TextValue $ws.Range("E32") "  +1.90%  "

Set-TextValue $ws.Range("D33") "154.97"
Set-TextValue $ws.Range("E33") "  +1.85%  "

Set-TextValue $ws.Range("D34") "4.53"
Set-TextValue $ws.Range("E34") "  +2.11%  "

Set-TextValue $ws.Range("D35") "6.01"
Set-TextValue $ws.Range("E35") "  +3.72%  "

Set-TextValue $ws.Range("D36") "26.92"
Set-TextValue $ws.Range("E36") "  +0.18%  "

Set-TextValue $ws.Range("E37") "  +4.53%  "

Set-TextValue $ws.Range("E38") "  +2.48%  "

Set-TextValue $ws.Range("D39") "3.107.30"
Set-TextValue $ws.Range("E39") "  +2.29%  "

Set-TextValue $ws.Range("D40") "3.91"
Set-TextValue $ws.Range("E40") "  +3.97%  "

Set-TextValue $ws.Range("D41") "36.49"
Set-TextValue $ws.Range("E41") "  +0.44%  "

Set-TextValue $ws.Range("D43") "0.654"
Set-TextValue $ws.Range("E43") "  +0.29%  "

Set-TextValue $ws.Range("D44") "2.260.65"
Set-TextValue $ws.Range("E44") "  +2.54%  "

Set-TextValue $ws.Range("D45") "0.0258"
Set-TextValue $ws.Range("E45") "  +8.12%  "

Set-TextValue $ws.Range("D46") "20.77"
Set-TextValue $ws.Range("E46") "  +6.31%  "

Set-TextValue $ws.Range("E47") "  +2.27%  "

Set-TextValue $ws.Range("D48") "0.942"
Set-TextValue $ws.Range("E48") "  +2.82%  "

Set-TextValue $ws.Range("D49") "5.92"
Set-TextValue $ws.Range("E49") "  +1.73%  "

Set-TextValue $ws.Range("D50") "0.739"
Set-TextValue $ws.Range("E50") "  +9.90%  "

Set-TextValue $ws.Range("D51") "256.77"
Set-TextValue $ws.Range("E51") "  +12.38%  "
